# beeswarm graphic work: add a log10(winningest) helper column to the
# "composite" sheet, widen the autofilter / filter-database range to cover
# it, and point the pivot caches back at their (re-synced) definitions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("composite")

# --- new column H: log_winningest -----------------------------------------
$ws.Range("H1").Value = "log_winningest"
$ws.Range("H2").Formula = "=LOG10(D2)"
$ws.Range("H3:H17").Formula = "=LOG10(D3)"

# --- move selection onto the newly-typed cell, like the author would have -
$ws.Activate()
$ws.Range("H2").Select()

# --- widen the autofilter so it covers the new column ----------------------
$ws.Range("A1:G1").AutoFilter()
$ws.Range("A1:G1").AutoFilter()

# --- widen the hidden _FilterDatabase defined name to match ---------------
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "composite!_FilterDatabase") {
        $n.RefersTo = "=composite!`$A`$1:`$G`$1"
    }
}
